# "Generate Report for Handback"
#
# For each localized-language sheet (zh-cn, de-de) this records that the
# previously handed-off file has now been handed back:
#   - Status (col B) moves from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - A new "Latest Target File" (col E) / "Latest Handback File" (col F)
#     pair of hyperlinks is recorded, mirroring the existing
#     Source File Name (col A) / Latest Handoff File (col C) links
#   - "Latest Handback DateTime" (col G) is stamped with the handback time

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# NOTE: this interpreter only binds POSITIONAL function parameters
# reliably, so every helper below is called positionally (no `-Name value`
# syntax).
function Update-LangSheet {
    param($SheetName, $HandbackDate, $MdUrl, $XlfUrl)

    $ws = $wb.Worksheets.Item($SheetName)

    # --- Status column: handed off -> handed back -------------------------
    $ws.Range("B2").Value = $statusText
    $ws.Range("B3").Value = $statusText

    # --- Latest Target File / Latest Handback File (new columns) ----------
    # Rows 2 and 3 mirror their own Source File Name (A) / Latest Handoff
    # File (C) hyperlinks.
    $mdDisplay = "9ef3cb4e-ac1a-4b6d-9780-cffdec1fe5b3.md"
    $xlfDisplay = "9ef3cb4e-ac1a-4b6d-9780-cffdec1fe5b3.a4ec450bc5649ca8d2b1692c9369a1628dc10bee." + $SheetName + ".xlf"

    $ws.Hyperlinks.Add($ws.Cells.Item(2, 5), $MdUrl, "", "", $mdDisplay) | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item(2, 6), $XlfUrl, "", "", $xlfDisplay) | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item(3, 5), $MdUrl, "", "", $mdDisplay) | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item(3, 6), $XlfUrl, "", "", $xlfDisplay) | Out-Null

    # --- Latest Handback DateTime ------------------------------------------
    $ws.Range("G2").Value = $HandbackDate
    $ws.Range("G3").Value = $HandbackDate
}

# zh-cn: handback recorded at 2016-03-09 08:44:09
Update-LangSheet "zh-cn" "2016-03-09 08:44:09" `
    "https://github.com/OpenLocalizationTest/oltest/blob/56de3b9fde643b69a374b2d5593db5246d397dd9/e2e/9ef3cb4e-ac1a-4b6d-9780-cffdec1fe5b3.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/169bb7c20c4041dc2d36065bbb88f3ef62ba86ee/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9ef3cb4e-ac1a-4b6d-9780-cffdec1fe5b3.a4ec450bc5649ca8d2b1692c9369a1628dc10bee.zh-cn.xlf"

# de-de: handback recorded at 2016-03-09 08:44:26
Update-LangSheet "de-de" "2016-03-09 08:44:26" `
    "https://github.com/OpenLocalizationTest/oltest/blob/56de3b9fde643b69a374b2d5593db5246d397dd9/e2e/9ef3cb4e-ac1a-4b6d-9780-cffdec1fe5b3.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/af7a1537a712888142c6461279eceaf7998ffcb1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9ef3cb4e-ac1a-4b6d-9780-cffdec1fe5b3.a4ec450bc5649ca8d2b1692c9369a1628dc10bee.de-de.xlf"

# --- Overview sheet: same status text rolls up for both languages ---------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

Write-Output "Handback report generated."
